# Regenerate the handback-status report hyperlinks/values for the new
# handback run (new source/target GUID file names + new timestamps).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name column hyperlinks
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewLinks = @(
    @{Addr="A2"; Url="https://github.com/OpenLocalizationTest/oltest/blob/64fb898706e8550345a8fdeb80a4de17d5f54a10/e2e/69ceb517-a00a-401b-81fe-83bf7f947dc1.md"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.md"},
    @{Addr="A3"; Url="https://github.com/OpenLocalizationTest/oltest/blob/64fb898706e8550345a8fdeb80a4de17d5f54a10/e2e/97022198-1c0d-4c28-be89-8765048b9914.md"; Disp="ffff0c4f82a7-3f0f-44ee-8b32-110fb29e5f60.md"}
)
$wsOverview.Hyperlinks.Delete()
foreach ($lnk in $overviewLinks) {
    $wsOverview.Hyperlinks.Add($wsOverview.Range($lnk.Addr), $lnk.Url, [Type]::Missing, [Type]::Missing, $lnk.Disp) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 13:13:33"
$wsZh.Range("H2").Value = "2016-03-23 13:13:56"
$wsZh.Range("E3").Value = "2016-03-23 13:13:33"
$wsZh.Range("H3").Value = "2016-03-23 13:13:56"

$zhLinks = @(
    @{Addr="A2"; Url="https://github.com/OpenLocalizationTest/oltest/blob/64fb898706e8550345a8fdeb80a4de17d5f54a10/e2e/69ceb517-a00a-401b-81fe-83bf7f947dc1.md"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.md"},
    @{Addr="D2"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8c1fbae451378bbb8bcc754b2bbd201e708b16ef/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/69ceb517-a00a-401b-81fe-83bf7f947dc1.806fa69c452a973fabe495c4f87272cbbade52fc.zh-cn.xlf"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.zh-cn.xlf"},
    @{Addr="F2"; Url="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/aad8c19bef06138c5cf509b2b791923cdb91123b/e2e/69ceb517-a00a-401b-81fe-83bf7f947dc1.md"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.md"},
    @{Addr="G2"; Url="https://github.com/OpenLocalizationTestOrg/olhandback/blob/73318a9fb81a4e32a0f88ac95e22377c41bb7c06/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/69ceb517-a00a-401b-81fe-83bf7f947dc1.806fa69c452a973fabe495c4f87272cbbade52fc.zh-cn.xlf"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.zh-cn.xlf"},
    @{Addr="A3"; Url="https://github.com/OpenLocalizationTest/oltest/blob/64fb898706e8550345a8fdeb80a4de17d5f54a10/e2e/97022198-1c0d-4c28-be89-8765048b9914.md"; Disp="ffff0c4f82a7-3f0f-44ee-8b32-110fb29e5f60.md"},
    @{Addr="D3"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8c1fbae451378bbb8bcc754b2bbd201e708b16ef/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/97022198-1c0d-4c28-be89-8765048b9914.906268ef18005d9894d02ef9be286407c6bfc634.zh-cn.xlf"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.zh-cn.xlf"},
    @{Addr="F3"; Url="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/aad8c19bef06138c5cf509b2b791923cdb91123b/e2e/97022198-1c0d-4c28-be89-8765048b9914.md"; Disp="ffff0c4f82a7-3f0f-44ee-8b32-110fb29e5f60.md"},
    @{Addr="G3"; Url="https://github.com/OpenLocalizationTestOrg/olhandback/blob/73318a9fb81a4e32a0f88ac95e22377c41bb7c06/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/97022198-1c0d-4c28-be89-8765048b9914.906268ef18005d9894d02ef9be286407c6bfc634.zh-cn.xlf"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.zh-cn.xlf"}
)
$wsZh.Hyperlinks.Delete()
foreach ($lnk in $zhLinks) {
    $wsZh.Hyperlinks.Add($wsZh.Range($lnk.Addr), $lnk.Url, [Type]::Missing, [Type]::Missing, $lnk.Disp) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 13:13:37"
$wsDe.Range("H2").Value = "2016-03-23 13:14:04"
$wsDe.Range("E3").Value = "2016-03-23 13:13:37"
$wsDe.Range("H3").Value = "2016-03-23 13:14:04"

$deLinks = @(
    @{Addr="A2"; Url="https://github.com/OpenLocalizationTest/oltest/blob/64fb898706e8550345a8fdeb80a4de17d5f54a10/e2e/69ceb517-a00a-401b-81fe-83bf7f947dc1.md"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.md"},
    @{Addr="D2"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/44b31fef91396b84f17dfb9898a2b33b53b0f783/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/69ceb517-a00a-401b-81fe-83bf7f947dc1.806fa69c452a973fabe495c4f87272cbbade52fc.de-de.xlf"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.de-de.xlf"},
    @{Addr="F2"; Url="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/64523117c1441b15fe31dd778c581ae1ee90f8f1/e2e/69ceb517-a00a-401b-81fe-83bf7f947dc1.md"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.md"},
    @{Addr="G2"; Url="https://github.com/OpenLocalizationTestOrg/olhandback/blob/2c2d98df55ccd4ee5a3d89e46ca6266b33ffb4b0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/69ceb517-a00a-401b-81fe-83bf7f947dc1.806fa69c452a973fabe495c4f87272cbbade52fc.de-de.xlf"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.de-de.xlf"},
    @{Addr="A3"; Url="https://github.com/OpenLocalizationTest/oltest/blob/64fb898706e8550345a8fdeb80a4de17d5f54a10/e2e/97022198-1c0d-4c28-be89-8765048b9914.md"; Disp="ffff0c4f82a7-3f0f-44ee-8b32-110fb29e5f60.md"},
    @{Addr="D3"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/44b31fef91396b84f17dfb9898a2b33b53b0f783/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/97022198-1c0d-4c28-be89-8765048b9914.906268ef18005d9894d02ef9be286407c6bfc634.de-de.xlf"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.de-de.xlf"},
    @{Addr="F3"; Url="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/64523117c1441b15fe31dd778c581ae1ee90f8f1/e2e/97022198-1c0d-4c28-be89-8765048b9914.md"; Disp="ffff0c4f82a7-3f0f-44ee-8b32-110fb29e5f60.md"},
    @{Addr="G3"; Url="https://github.com/OpenLocalizationTestOrg/olhandback/blob/2c2d98df55ccd4ee5a3d89e46ca6266b33ffb4b0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/97022198-1c0d-4c28-be89-8765048b9914.906268ef18005d9894d02ef9be286407c6bfc634.de-de.xlf"; Disp="1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.de-de.xlf"}
)
$wsDe.Hyperlinks.Delete()
foreach ($lnk in $deLinks) {
    $wsDe.Hyperlinks.Add($wsDe.Range($lnk.Addr), $lnk.Url, [Type]::Missing, [Type]::Missing, $lnk.Disp) | Out-Null
}
